$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that gets bumped by one day
# (45180 -> 45181) for every data row (rows 2 through 250).
$ws.Range("C2:C250").Value = 45181
